$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.999.64"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.305.39"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.12"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.64%  "

$ws.Range("E7").Value = "  -1.77%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -0.70%  "

$ws.Range("E10").Value = "  -0.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.74"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.15%  "

$ws.Range("E14").Value = "  +2.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.665.57"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.308.85"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("E17").Value = "  +1.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.881.92"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("E19").Value = "  -1.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.66%  "

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.73"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.90"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.46%  "

$ws.Range("E24").Value = "  +1.81%  "

$ws.Range("E25").Value = "  +0.17%  "

$ws.Range("E26").Value = "  +0.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.91"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +8.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.18"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.06"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.83"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.20"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.29%  "

$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.48"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.99%  "

$ws.Range("E36").Value = "  -1.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0690"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.82%  "

$ws.Range("E38").Value = "  -0.40%  "

$ws.Range("E39").Value = "  +0.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.97%  "

$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.996.87"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.08%  "

$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.35"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.71%  "

$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.84"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.66%  "

$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.531.55"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.45"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("E50").Value = "  -3.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.13%  "
